$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "30.555.55"
Set-TextValue $ws.Range("E2") "  -0.19%  "
Set-TextValue $ws.Range("D3") "1.913.70"
Set-TextValue $ws.Range("E3") "  -0.49%  "
Set-TextValue $ws.Range("D4") "0.9998"
Set-TextValue $ws.Range("E4") "  -0.05%  "
Set-TextValue $ws.Range("D5") "244.56"
Set-TextValue $ws.Range("E5") "  -0.75%  "
Set-TextValue $ws.Range("D6") "0.9999"
Set-TextValue $ws.Range("E6") "  -0.04%  "
Set-TextValue $ws.Range("D7") "0.4850"
Set-TextValue $ws.Range("E7") "  +2.30%  "
Set-TextValue $ws.Range("E8") "  +0.28%  "
Set-TextValue $ws.Range("D9") "0.06805"
Set-TextValue $ws.Range("E9") "  -0.31%  "
Set-TextValue $ws.Range("D10") "111.15"
Set-TextValue $ws.Range("E10") "  +5.54%  "
Set-TextValue $ws.Range("D11") "19.37"
Set-TextValue $ws.Range("E11") "  +5.60%  "
Set-TextValue $ws.Range("D12") "1.917.53"
Set-TextValue $ws.Range("E12") "  -0.31%  "
Set-TextValue $ws.Range("D13") "0.07571"
Set-TextValue $ws.Range("E13") "  -1.59%  "
Set-TextValue $ws.Range("D14") "5.381"
Set-TextValue $ws.Range("E14") "  +0.73%  "
Set-TextValue $ws.Range("D15") "0.6710"
Set-TextValue $ws.Range("E15") "  +0.55%  "
Set-TextValue $ws.Range("D16") "296.09"
Set-TextValue $ws.Range("E16") "  +1.43%  "
Set-TextValue $ws.Range("D17") "30.555.12"
Set-TextValue $ws.Range("E17") "  -0.18%  "
Set-TextValue $ws.Range("D18") "13.02"
Set-TextValue $ws.Range("E18") "  +0.56%  "
Set-TextValue $ws.Range("D19") "0.000007595"
Set-TextValue $ws.Range("E19") "  -0.14%  "
Set-TextValue $ws.Range("D20") "1.0000"
Set-TextValue $ws.Range("E20") "  -0.01%  "
Set-TextValue $ws.Range("D21") "5.532"
Set-TextValue $ws.Range("E21") "  -0.27%  "
Set-TextValue $ws.Range("D22") "2.163.54"
Set-TextValue $ws.Range("E22") "  -0.21%  "
Set-TextValue $ws.Range("D23") "0.9998"
Set-TextValue $ws.Range("E23") "  -0.11%  "
Set-TextValue $ws.Range("D24") "6.451"
Set-TextValue $ws.Range("E24") "  +0.37%  "
Set-TextValue $ws.Range("D25") "9.472"
Set-TextValue $ws.Range("E25") "  +0.31%  "
Set-TextValue $ws.Range("D26") "165.96"
Set-TextValue $ws.Range("E26") "  -0.94%  "
Set-TextValue $ws.Range("D27") "20.31"
Set-TextValue $ws.Range("E27") "  -3.58%  "
Set-TextValue $ws.Range("D28") "2.078"
Set-TextValue $ws.Range("E28") "  -1.64%  "
Set-TextValue $ws.Range("E29") "  -0.36%  "
Set-TextValue $ws.Range("D30") "1.435"
Set-TextValue $ws.Range("E30") "  +2.65%  "
Set-TextValue $ws.Range("D31") "4.149"
Set-TextValue $ws.Range("E31") "  -0.70%  "
Set-TextValue $ws.Range("D32") "4.054"
Set-TextValue $ws.Range("E32") "  +0.03%  "
Set-TextValue $ws.Range("E33") "  -0.83%  "
Set-TextValue $ws.Range("D34") "0.7349"
Set-TextValue $ws.Range("E34") "  -0.26%  "
Set-TextValue $ws.Range("D35") "1.141"
Set-TextValue $ws.Range("E35") "  -0.25%  "
Set-TextValue $ws.Range("D36") "0.9993"
Set-TextValue $ws.Range("E36") "  -0.01%  "
Set-TextValue $ws.Range("E37") "  -1.15%  "
Set-TextValue $ws.Range("D38") "2.718"
Set-TextValue $ws.Range("E38") "  -0.78%  "
Set-TextValue $ws.Range("E39") "  -0.03%  "
Set-TextValue $ws.Range("D40") "2.020"
Set-TextValue $ws.Range("E40") "  -1.59%  "
Set-TextValue $ws.Range("D41") "109.17"
Set-TextValue $ws.Range("E41") "  -1.85%  "
Set-TextValue $ws.Range("D42") "0.4447"
Set-TextValue $ws.Range("E42") "  +2.02%  "
Set-TextValue $ws.Range("D43") "0.8674"
Set-TextValue $ws.Range("E43") "  -0.64%  "
Set-TextValue $ws.Range("D44") "5.815"
Set-TextValue $ws.Range("E44") "  -1.75%  "
Set-TextValue $ws.Range("D45") "0.9998"
Set-TextValue $ws.Range("E45") "  -0.04%  "
Set-TextValue $ws.Range("D46") "69.59"
Set-TextValue $ws.Range("E46") "  +2.55%  "
Set-TextValue $ws.Range("D47") "7.209"
Set-TextValue $ws.Range("E47") "  -0.76%  "
Set-TextValue $ws.Range("D48") "48.47"
Set-TextValue $ws.Range("E48") "  +0.56%  "
Set-TextValue $ws.Range("D49") "9.253"
Set-TextValue $ws.Range("E49") "  -0.65%  "
Set-TextValue $ws.Range("D50") "0.1226"
Set-TextValue $ws.Range("E50") "  -1.35%  "
Set-TextValue $ws.Range("D51") "0.2509"
Set-TextValue $ws.Range("E51") "  +0.53%  "
